# Update the salted output data (column B) with the newly generated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(-278, 200, 323, 388, -313, 418, -244, 12, 98, -99, 802, 739, 838, 1166, 636, 1230, 841, 1541, 1261, 1851, 1505, 1700, 2388, 2663, 2178, 2741, 2374, 3427, 3391, 3838, 4183, 4193, 4411, 4062, 4341, 4813, 5575, 5985, 6123, 6532, 6202, 6822, 7451, 7498, 7547, 8814, 8388, 8609, 9709, 9745, 10542, 11028, 10900, 11795, 11538, 12002, 13047, 12879, 13287, 14644, 14555, 15426, 15457, 16330, 17181, 16766, 17312, 18639, 18516, 19868, 20394, 20968, 20816, 21902, 22008, 23101, 23298, 23622, 25108, 24993, 26448, 26226, 27019, 28103, 29099, 29567, 30377, 30918, 31702, 32569, 32740, 33180, 34188, 35451, 36126, 36920, 37466, 37752, 38639, 39504)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Match the author's final selection state in the saved workbook.
$ws.Range("L14").Select()
